$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 600
$ws.Range("B3").Value = 350
$ws.Range("B4").Value = 97
$ws.Range("B6").Value = 300
